$wb = $excel.ActiveWorkbook

# Update regression coefficients/std-errors on the six "_Males"/"_Females" summary sheets

$ws = $wb.Worksheets.Item("UK_DWB_MCS2_Males")
$ws.Range("B2").Value = [double]"-1.23056636431974"
$ws.Range("C2").Value = [double]"0.103445250024589"
$ws.Range("B3").Value = [double]"1.4441139064634101"
$ws.Range("D3").Value = [double]"6.6087181550199103E-2"
$ws.Range("B4").Value = [double]"-0.58348744071858405"
$ws.Range("E4").Value = [double]"3.2081096346227099E-2"
$ws.Range("B5").Value = [double]"-0.25042214934580498"
$ws.Range("F5").Value = [double]"4.6389313542080403E-2"
$ws.Range("B6").Value = [double]"-0.86343025909766202"
$ws.Range("G6").Value = [double]"6.3548920856129704E-2"
$ws.Range("B7").Value = [double]"-1.0891121745363399"
$ws.Range("H7").Value = [double]"6.8177272558643906E-2"
$ws.Range("B8").Value = [double]"2.5702068477583499E-2"
$ws.Range("I8").Value = [double]"1.5105874110923199E-2"
$ws.Range("B9").Value = [double]"-0.16470758086589499"
$ws.Range("J9").Value = [double]"1.1374621686648E-2"

$ws = $wb.Worksheets.Item("UK_DWB_MCS2_Females")
$ws.Range("B2").Value = [double]"-1.4228190123050499"
$ws.Range("C2").Value = [double]"0.209206184901087"
$ws.Range("B3").Value = [double]"1.0307197943903601"
$ws.Range("D3").Value = [double]"9.4849755363602295E-2"
$ws.Range("B4").Value = [double]"-1.3672336758516399"
$ws.Range("E4").Value = [double]"0.117003284367741"
$ws.Range("B5").Value = [double]"-0.46015958776090499"
$ws.Range("F5").Value = [double]"7.2503426898626294E-2"
$ws.Range("B6").Value = [double]"-0.70629664050724295"
$ws.Range("G6").Value = [double]"9.82289801979849E-2"
$ws.Range("B7").Value = [double]"-1.2474728694726001"
$ws.Range("H7").Value = [double]"0.12303306255009799"
$ws.Range("B8").Value = [double]"-9.5896950946945704E-3"
$ws.Range("I8").Value = [double]"1.57236914292413E-2"
$ws.Range("B9").Value = [double]"-0.20401663530286501"
$ws.Range("J9").Value = [double]"1.33149079995976E-2"

$ws = $wb.Worksheets.Item("UK_DWB_PCS2_Males")
$ws.Range("B2").Value = [double]"0.14473670251260101"
$ws.Range("C2").Value = [double]"6.1273226234816798E-2"
$ws.Range("B3").Value = [double]"0.72398980017675996"
$ws.Range("D3").Value = [double]"4.2677456089765899E-2"
$ws.Range("B4").Value = [double]"-0.40805728765892602"
$ws.Range("E4").Value = [double]"2.5804704114682799E-2"
$ws.Range("B5").Value = [double]"-0.64786512215603698"
$ws.Range("F5").Value = [double]"3.7210638363925901E-2"
$ws.Range("B6").Value = [double]"-0.52790238319098104"
$ws.Range("G6").Value = [double]"5.2968602686544401E-2"
$ws.Range("B7").Value = [double]"-0.71118265552130999"
$ws.Range("H7").Value = [double]"5.4453392923858999E-2"
$ws.Range("B8").Value = [double]"3.5313846726256898E-2"
$ws.Range("I8").Value = [double]"1.0665778176003801E-2"
$ws.Range("B9").Value = [double]"0.14230834922964999"
$ws.Range("J9").Value = [double]"8.1244931188861596E-3"

$ws = $wb.Worksheets.Item("UK_DWB_PCS2_Females")
$ws.Range("B2").Value = [double]"0.96700780202930603"
$ws.Range("C2").Value = [double]"0.11319861709536"
$ws.Range("B3").Value = [double]"-8.7236005864526206E-2"
$ws.Range("D3").Value = [double]"7.3896310888152197E-2"
$ws.Range("B4").Value = [double]"-0.32979624167820298"
$ws.Range("E4").Value = [double]"7.8297133540976296E-2"
$ws.Range("B5").Value = [double]"-0.72965074119122597"
$ws.Range("F5").Value = [double]"5.8308252125873698E-2"
$ws.Range("B6").Value = [double]"-0.70814708182476604"
$ws.Range("G6").Value = [double]"7.6068346940108E-2"
$ws.Range("B7").Value = [double]"-1.0288942592311801"
$ws.Range("H7").Value = [double]"8.2842310652720505E-2"
$ws.Range("B8").Value = [double]"-9.5008964591471401E-2"
$ws.Range("I8").Value = [double]"9.8134737056849008E-3"
$ws.Range("B9").Value = [double]"6.1434528696179998E-2"
$ws.Range("J9").Value = [double]"8.7577694998627904E-3"

$ws = $wb.Worksheets.Item("UK_DLS2_Males")
$ws.Range("B2").Value = [double]"-9.9638023612954804E-2"
$ws.Range("C2").Value = [double]"2.1877917000088401E-3"
$ws.Range("B3").Value = [double]"6.12046933064666E-2"
$ws.Range("D3").Value = [double]"1.7090782019972899E-3"
$ws.Range("B4").Value = [double]"-7.9217973407544703E-2"
$ws.Range("E4").Value = [double]"7.6124667237600896E-4"
$ws.Range("B5").Value = [double]"-0.13470996687658601"
$ws.Range("F5").Value = [double]"1.2411083062172201E-3"
$ws.Range("B6").Value = [double]"-0.115184567399258"
$ws.Range("G6").Value = [double]"1.44347865994554E-3"
$ws.Range("B7").Value = [double]"-0.153507417191036"
$ws.Range("H7").Value = [double]"1.50291263094377E-3"
$ws.Range("B8").Value = [double]"-1.21940132240634E-2"
$ws.Range("I8").Value = [double]"3.6956120195602002E-4"
$ws.Range("B9").Value = [double]"-1.61707100239993E-2"
$ws.Range("J9").Value = [double]"2.5383801812351998E-4"

$ws = $wb.Worksheets.Item("UK_DLS2_Females")
$ws.Range("B2").Value = [double]"-0.241099051541424"
$ws.Range("C2").Value = [double]"4.4857959823730101E-3"
$ws.Range("B3").Value = [double]"7.2105821167758502E-2"
$ws.Range("D3").Value = [double]"3.1340383658325301E-3"
$ws.Range("B4").Value = [double]"-0.33578763811511397"
$ws.Range("E4").Value = [double]"2.6440006564389502E-3"
$ws.Range("B5").Value = [double]"-0.141027262742555"
$ws.Range("F5").Value = [double]"1.91969597551501E-3"
$ws.Range("B6").Value = [double]"-0.185338761984772"
$ws.Range("G6").Value = [double]"2.2383139962641502E-3"
$ws.Range("B7").Value = [double]"-0.172122215639476"
$ws.Range("H7").Value = [double]"2.3839546743497902E-3"
$ws.Range("B8").Value = [double]"7.0854464583095999E-3"
$ws.Range("I8").Value = [double]"3.9984848987148697E-4"
$ws.Range("B9").Value = [double]"-4.36990950739094E-2"
$ws.Range("J9").Value = [double]"3.6007712287628498E-4"

# Leave the workbook scrolled on UK_DWB_MCS1 (A25) with UK_DWB_MCS2_Males as the active tab
$wsFirst = $wb.Worksheets.Item("UK_DWB_MCS1")
$wsFirst.Activate()
$wsFirst.Range("A25").Select()
$wb.Worksheets.Item("UK_DWB_MCS2_Males").Activate()
